# PSP_Sheet.xlsx update: add new Time Recording Log entries (11/17, 11/21,
# 11/27, 12/02, 12/03) to the "김수인" sheet (Sheet2) and refresh the
# on-screen scroll/selection state of the two personal-log sheets, matching
# commit "update report, PSP, django project file".

$wb = $excel.ActiveWorkbook
$sheetJoonki = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# Row 26 - 11월 17일 / Django Template
# ---------------------------------------------------------------------
$ws.Range("A25:F25").Copy($ws.Range("A26:F26"))
$ws.Rows(26).RowHeight = 13.8
$ws.Range("A26").Value = "11월 17일"
$rA26 = $ws.Range("A26")
$rA26.Characters(3, 4).Font.Name = "Arial"
$rA26.Characters(3, 4).Font.Size = 10
$rA26.Characters(7, 1).Font.Name = "돋움"
$rA26.Characters(7, 1).Font.Size = 10
$ws.Range("B26").Value = 0.083333333333333329
$ws.Range("C26").Value = 0.29166666666666669
$ws.Range("D26").Value = 120
$ws.Range("E26").Value = 300
$ws.Range("F26").Value = "Django Template 개발 소스코드 Bootstrap 실습"
$rF26 = $ws.Range("F26")
$rF26.Characters(19, 18).Font.Name = "돋움"
$rF26.Characters(19, 18).Font.Size = 10

# ---------------------------------------------------------------------
# Row 27 - 11월 21일 / Diagram feedback
# ---------------------------------------------------------------------
$ws.Range("A25:F25").Copy($ws.Range("A27:F27"))
$ws.Rows(27).RowHeight = 13.8
$ws.Range("A27").Value = "11월 21일"
$rA27 = $ws.Range("A27")
$rA27.Characters(3, 4).Font.Name = "Arial"
$rA27.Characters(3, 4).Font.Size = 10
$rA27.Characters(7, 1).Font.Name = "돋움"
$rA27.Characters(7, 1).Font.Size = 10
$ws.Range("B27").Value = 0.5
$ws.Range("C27").Value = 0.54166666666666663
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 60
$ws.Range("F27").Value = "Diagram 피드백 참여"
$rF27 = $ws.Range("F27")
$rF27.Characters(12, 3).Font.Name = "돋움"
$rF27.Characters(12, 3).Font.Size = 10

# ---------------------------------------------------------------------
# Row 28 - 11월 27일 / Django code, result page capture
# ---------------------------------------------------------------------
$ws.Range("A25:F25").Copy($ws.Range("A28:F28"))
$ws.Rows(28).RowHeight = 15.6
$ws.Range("A28").Value = "11월 27일"
$rA28 = $ws.Range("A28")
$rA28.Characters(3, 1).Font.Name = "맑은 고딕"
$rA28.Characters(3, 1).Font.Size = 10
$rA28.Characters(4, 3).Font.Name = "Arial"
$rA28.Characters(4, 3).Font.Size = 10
$rA28.Characters(7, 1).Font.Name = "돋움"
$rA28.Characters(7, 1).Font.Size = 10
$ws.Range("B28").Value = 0.625
$ws.Range("C28").Value = 0.75
$ws.Range("D28").Value = 30
$ws.Range("E28").Value = 150
$ws.Range("F28").Value = "Django code, result page 캡쳐, diagram upload"
$rF28 = $ws.Range("F28")
$rF28.Characters(28, 16).Font.Name = "돋움"
$rF28.Characters(28, 16).Font.Size = 10

# ---------------------------------------------------------------------
# Row 29 - 12월 02일 / Django Chart
# ---------------------------------------------------------------------
$ws.Range("A25:F25").Copy($ws.Range("A29:F29"))
$ws.Rows(29).RowHeight = 13.8
$ws.Range("A29").Value = "12월 02일"
$rA29 = $ws.Range("A29")
$rA29.Characters(4, 4).Font.Name = "돋움"
$rA29.Characters(4, 4).Font.Size = 10
$ws.Range("B29").Value = 0.91666666666666663
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 20
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = "Django Chart 구현"

# ---------------------------------------------------------------------
# Row 30 - 12월 03일 / Django application, report & PSP sheet update
# ---------------------------------------------------------------------
$ws.Range("A25:F25").Copy($ws.Range("A30:F30"))
$ws.Rows(30).RowHeight = 13.8
$ws.Range("A30").Value = "12월 03일"
$rA30 = $ws.Range("A30")
$rA30.Characters(3, 4).Font.Name = "Arial"
$rA30.Characters(3, 4).Font.Size = 10
$rA30.Characters(7, 1).Font.Name = "돋움"
$rA30.Characters(7, 1).Font.Size = 10
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0.25
$ws.Range("D30").Value = 60
$ws.Range("E30").Value = 300
$ws.Range("F30").Value = "Django application 구현, 보고서 수정, PSP sheet 수정 "
$rF30 = $ws.Range("F30")
$rF30.Characters(22, 23).Font.Name = "돋움"
$rF30.Characters(22, 23).Font.Size = 10

# ---------------------------------------------------------------------
# Refresh on-screen scroll/selection state for both personal-log sheets.
# "이준기" (Sheet1) scrolled down to row 19; the selected cell (J21) is
# unchanged. "김수인" (Sheet2) scrolled to row 16 and the new last entry
# (F30) is now selected. Sheet2 stays the active tab, as before.
# ---------------------------------------------------------------------
$sheetJoonki.Activate()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F30").Select() | Out-Null
